# Insert a new data row before the current row 926, shifting all rows
# from 926-1042 down to 927-1043, then populate the newly inserted row
# 926 with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a whole new row at position 926 (existing row 926 and below move down by one)
$ws.Rows("926:926").Insert()

# Populate the new row 926 with the new record
$ws.Cells.Item(926, 1).Value = 10                               # A - Mercado ID
$ws.Cells.Item(926, 2).Value = "Vega Modelo de Temuco"          # B - Mercado
$ws.Cells.Item(926, 3).Value = "La Araucanía"                   # C - Región
$ws.Cells.Item(926, 4).Value = 45077                            # D - Fecha
$ws.Cells.Item(926, 5).Value = 9                                # E - Codreg
$ws.Cells.Item(926, 6).Value = 100112006                        # F - Categoría ID
$ws.Cells.Item(926, 7).Value = "Repollo"                        # G - Categoría
$ws.Cells.Item(926, 8).Value = "Crespo record"                  # H - Variedad
$ws.Cells.Item(926, 9).Value = "Primera"                        # I - Calidad
$ws.Cells.Item(926, 10).Value = 1250                            # J - Volumen
$ws.Cells.Item(926, 11).Value = 1200                            # K - Precio mínimo
$ws.Cells.Item(926, 12).Value = 1200                            # L - Precio máximo
$ws.Cells.Item(926, 13).Value = 1200                            # M - Precio promedio ponderado
$ws.Cells.Item(926, 14).Value = "$/unidad"                      # N - Unidad de comercialización
$ws.Cells.Item(926, 15).Value = "Región del Maule"               # O - Origen
$ws.Cells.Item(926, 16).Value = 1200                            # P - Precio $/Kg
$ws.Cells.Item(926, 17).Value = 1                               # Q - Kg o Unidades
$ws.Cells.Item(926, 18).Value = "Hortaliza"                     # R - Clasificación

# Restore the date number format on the newly inserted date cell
$ws.Cells.Item(926, 4).NumberFormat = $ws.Cells.Item(927, 4).NumberFormat
